$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update results for previously-pending rows 4, 10 and 12 ---
$ws.Range("G4").Value = "Fallo"
$ws.Range("H4").Value = -1

$ws.Range("G10").Value = "Fallo"
$ws.Range("H10").Value = -1

$ws.Range("G12").Value = "Fallo"
$ws.Range("H12").Value = -1

# --- Append new tracked events as rows 15 and 16 ---
$ws.Range("A15").Value = 14714061
$ws.Range("B15").Value = "'2025-09-19"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "Alejandro Tabilo"
$ws.Range("D15").Value = "Luciano Darderi"
$ws.Range("E15").Value = "Gana Luciano Darderi"
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = "'"
$ws.Range("G15").Style = "Normal"
$ws.Range("H15").Value = "'"
$ws.Range("H15").Style = "Normal"

$ws.Range("A16").Value = 14713398
$ws.Range("B16").Value = "'2025-09-18"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "Millen Hurrion"
$ws.Range("D16").Value = "Timofei Derepasko"
$ws.Range("E16").Value = "Gana Timofei Derepasko"
$ws.Range("F16").Value = 1.83
$ws.Range("G16").Value = "'"
$ws.Range("G16").Style = "Normal"
$ws.Range("H16").Value = "'"
$ws.Range("H16").Style = "Normal"
